# Clarify that an MLCC (ceramic) capacitor is acceptable in place of the
# tantalum part that was previously called out as the primary choice for
# the 4.7u / C9,C13,C17,C21 BOM line.
#
# Swap the primary Mouser part number to a ceramic (MLCC) part, and update
# the alternate note to reference the tantalum part (formerly primary) as
# the alternate option.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G17").Value = "GRM188R61E475KE11D"
$ws.Range("H17").Value = "Or tantalum egAVX F981C475MMA"
